# Atualiza a coluna C (Ano): de "31/12/AAAA" para "01/01/AAAA",
# mantendo o mesmo ano em cada linha (programação das figuras do grupo 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 91 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -match '^31/12/(\d{4})$') {
        $year = $matches[1]
        # Mantém a célula como texto (evita a conversão automática para
        # número de série de data que o Excel faria ao reconhecer o padrão),
        # depois limpa a formatação temporária para preservar o estilo original.
        $cell.NumberFormat = "@"
        $cell.Value = "01/01/$year"
        $cell.ClearFormats()
    }
}
